# Update the "cryptos" price/volume table (columns D = Price, E = Volume(1h))
# to the latest scraped snapshot. Only the cells whose values actually changed
# are touched; everything else (A/B/C columns, headers, row 1) is left alone.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.151.64"

$ws.Range("D3").Value = "1.831.34"
$ws.Range("E3").Value = "  -0.10%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.59"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6604"
$ws.Range("E6").Value = "  -1.59%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07401"
$ws.Range("E8").Value = "  -0.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2928"
$ws.Range("E9").Value = "  -1.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.89"
$ws.Range("E10").Value = "  +0.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07761"
$ws.Range("E11").Value = "  +1.44%  "

$ws.Range("D12").Value = "1.830.60"
$ws.Range("E12").Value = "  -0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.995"
$ws.Range("E13").Value = "  -0.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6659"
$ws.Range("E14").Value = "  -1.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.84"
$ws.Range("E15").Value = "  -4.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.116"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008428"
$ws.Range("E17").Value = "  +2.49%  "

$ws.Range("D18").Value = "29.158.42"
$ws.Range("E18").Value = "  +0.33%  "

$ws.Range("D19").Value = "2.079.10"
$ws.Range("E19").Value = "  +0.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.07"
$ws.Range("E20").Value = "  +0.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.43"
$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.127"
$ws.Range("E23").Value = "  -2.46%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9998"
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.61"
$ws.Range("E25").Value = "  -0.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.598"
$ws.Range("E26").Value = "  -0.90%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1389"
$ws.Range("E27").Value = "  -2.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.92"
$ws.Range("E28").Value = "  -0.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.516"
$ws.Range("E29").Value = "  +0.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.113"
$ws.Range("E30").Value = "  -2.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.038"
$ws.Range("E31").Value = "  -1.87%  "

$ws.Range("E32").Value = "  -0.75%  "

$ws.Range("E33").Value = "  -2.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.865"
$ws.Range("E34").Value = "  +0.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7403"
$ws.Range("E35").Value = "  -1.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.141"
$ws.Range("E36").Value = "  +1.57%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.652"
$ws.Range("E37").Value = "  -1.13%  "

$ws.Range("D38").Value = "1.300.18"
$ws.Range("E38").Value = "  +0.46%  "

$ws.Range("E39").Value = "  -0.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.732"
$ws.Range("E40").Value = "  +0.95%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9207"
$ws.Range("E41").Value = "  -1.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.943"
$ws.Range("E42").Value = "  -2.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.08487"
$ws.Range("E43").Value = "  +2.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.18"
$ws.Range("E45").Value = "  -1.75%  "

$ws.Range("D46").Value = "1.974.44"
$ws.Range("E46").Value = "  +0.10%  "

$ws.Range("E47").Value = "  -0.66%  "

$ws.Range("E48").Value = "  -10.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.749"
$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.11"
$ws.Range("E50").Value = "  -0.59%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05844"
$ws.Range("E51").Value = "  -1.33%  "
